$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 42651.599502314813
$ws.Range("B2").Value = 14
$ws.Range("D2").Value = 46
$ws.Range("E2").Value = 6821
$ws.Range("F2").Value = 381
$ws.Range("G2").Value = 57
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 94
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 42446
$ws.Range("L2").Value = 73
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = 51
$ws.Range("O2").Value = 3
$ws.Range("Q2").Value = 46.242130528922125
$ws.Range("W2").Value = 1
